$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# The status text changed everywhere it is used (shared string), from
# "Ready for handoff" to "Handed back: in sync with en-US".
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"

# Updated handback datetimes.
$wsZhCn.Range("K2").Value = "2016-08-14 01:06:51"
$wsDeDe.Range("K2").Value = "2016-08-14 01:07:03"

# Error detail cleared now that the handback is in sync.
$wsZhCn.Range("P2").Value = ""
$wsDeDe.Range("P2").Value = ""

# Widen columns to accommodate the longer status text.
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("P1").ColumnWidth = 13.7470528738839

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("P1").ColumnWidth = 13.7470528738839
